$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.201.17"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.110.32"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'579.63"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'173.38"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "'6.55"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'0.477"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D13").Value = "'36.77"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "3.626.14"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "67.144.24"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "3.107.91"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'16.54"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "'491.46"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'0.702"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'7.84"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").Value = "'83.93"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'13.09"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "'10.60"
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("D28").Value = "'7.91"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'28.30"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").Value = "'0.114"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "0.0₃0932"
$ws.Range("E33").Value = "  -7.12%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'5.83"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").Value = "'0.972"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").Value = "'47.11"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "'2.04"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").Value = "'0.308"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("D42").Value = "'386.79"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "2.800.71"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "'2.56"
$ws.Range("E44").Value = "  -8.59%  "
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "'135.03"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'24.98"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "'2.20"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'6.70"
$ws.Range("E51").Value = "  -2.68%  "
